# Upload new version with timestamp
# Adds two new inventory rows ("ماسك جلسات اطفال" and "محلول ملح") right
# before the existing "مناديل سولو سحب" row (which shifts down one row),
# updates the running total and refreshes the footer timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

# ---------------------------------------------------------------------
# 1. Make room: insert two fresh rows right above the current total row
#    (row 39). This pushes the old total row (39 -> 41) and the footer
#    row (40 -> 42) down, Excel-style, carrying their merges with them.
# ---------------------------------------------------------------------
$ws.Rows.Item(39).Insert()
$ws.Rows.Item(39).Insert()

# ---------------------------------------------------------------------
# 2. Give the two new rows (39 & 40) the same per-cell formatting as the
#    existing data row 38 (borders/fonts/number formats/fills), and match
#    its row heights.
# ---------------------------------------------------------------------
foreach ($col in $cols) {
    $ws.Range($col + "38").Copy()
    $ws.Range($col + "39").PasteSpecial($xlPasteFormats)
    $ws.Range($col + "38").Copy()
    $ws.Range($col + "40").PasteSpecial($xlPasteFormats)
}
$ws.Application.CutCopyMode = $false

$ws.Rows.Item(39).RowHeight = 25.5
$ws.Rows.Item(40).RowHeight = 24.75

# Recreate the merges for the two new data rows, matching the pattern
# used by every other item row (e.g. row 38).
$ws.Range("A39:B39").Merge()
$ws.Range("C39:G39").Merge()
$ws.Range("H39:K39").Merge()
$ws.Range("L39:M39").Merge()
$ws.Range("N39:O39").Merge()

$ws.Range("A40:B40").Merge()
$ws.Range("C40:G40").Merge()
$ws.Range("H40:K40").Merge()
$ws.Range("L40:M40").Merge()
$ws.Range("N40:O40").Merge()

# ---------------------------------------------------------------------
# 3. Helper to write a "numeric looking" text value (e.g. "20.0000",
#    "16:0", "0:0") into a cell without Excel silently turning it into a
#    real number / time. We briefly switch the cell to text format, set
#    the value, then restore its normal display format.
# ---------------------------------------------------------------------
function Set-TextValue($addr, $numFmt, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = $numFmt
}

# ---------------------------------------------------------------------
# 4. Row 38 becomes the new item "ماسك جلسات اطفال".
# ---------------------------------------------------------------------
Set-TextValue "C38" "@" "ماسك جلسات اطفال"
Set-TextValue "H38" "@" "0:0"
Set-TextValue "N38" "@" "20.00"
Set-TextValue "P38" "0.00" "20.0000"

# ---------------------------------------------------------------------
# 5. Row 39 becomes the new item "محلول ملح".
# ---------------------------------------------------------------------
$ws.Range("A39").Value = 33
Set-TextValue "C39" "@" "محلول ملح"
Set-TextValue "H39" "@" "19:0"
Set-TextValue "L39" "#,##0.##;`"[`"#,##0.##`"]`";0" "0"
Set-TextValue "N39" "@" "24.00"
Set-TextValue "P39" "0.00" "24.0000"
Set-TextValue "Q39" "@" "1:0"

# ---------------------------------------------------------------------
# 6. Row 40 is the original "مناديل سولو سحب" item, now shifted down.
# ---------------------------------------------------------------------
$ws.Range("A40").Value = 34
Set-TextValue "C40" "@" "مناديل سولو سحب"
Set-TextValue "H40" "@" "16:0"
Set-TextValue "L40" "#,##0.##;`"[`"#,##0.##`"]`";0" "0"
Set-TextValue "N40" "@" "45.00"
Set-TextValue "P40" "0.00" "45.0000"
Set-TextValue "Q40" "@" "1:0"

# ---------------------------------------------------------------------
# 7. Row 41 (old row 39) is the running-total row: bump the total to
#    reflect the two newly added items.
# ---------------------------------------------------------------------
$ws.Range("P41").Value = 1521.9649999999999

# ---------------------------------------------------------------------
# 8. Row 42 (old row 40) is the footer: refresh the generated timestamp.
# ---------------------------------------------------------------------
$ws.Range("A42").Value = "Monday, 18 August, 2025 1:36 PM"

$wb.Save()
